$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep Price/Volume columns as plain text so values like "1.000" are not
# auto-converted to numbers by Excel when assigned.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.952.66'
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = '1.894.06'
$ws.Range("E3").Value = '  +0.14%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '0.7728'
$ws.Range("E5").Value = '  -1.99%  '
$ws.Range("D6").Value = '244.40'
$ws.Range("E6").Value = '  +0.66%  '
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").Value = '0.3137'
$ws.Range("E8").Value = '  -0.18%  '
$ws.Range("D9").Value = '25.69'
$ws.Range("E9").Value = '  +1.99%  '
$ws.Range("D10").Value = '0.07316'
$ws.Range("E10").Value = '  +3.81%  '
$ws.Range("D11").Value = '0.08047'
$ws.Range("E11").Value = '  -0.35%  '
$ws.Range("D12").Value = '0.7729'
$ws.Range("E12").Value = '  +1.20%  '
$ws.Range("D13").Value = '5.509'
$ws.Range("E13").Value = '  +2.92%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.878.33'
$ws.Range("E14").Value = '  -0.16%  '
$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").Value = '94.79'
$ws.Range("E15").Value = '  +3.00%  '
$ws.Range("D16").Value = '6.203'
$ws.Range("E16").Value = '  +4.01%  '
$ws.Range("D17").Value = '29.850.54'
$ws.Range("D18").Value = '14.03'
$ws.Range("E18").Value = '  +1.80%  '
$ws.Range("D19").Value = '247.53'
$ws.Range("E19").Value = '  +2.13%  '
$ws.Range("D20").Value = '0.000007847'
$ws.Range("E20").Value = '  +2.27%  '
$ws.Range("D21").Value = '8.179'
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("D22").Value = '0.9985'
$ws.Range("E22").Value = '  -0.31%  '
$ws.Range("D23").Value = '2.084.31'
$ws.Range("E23").Value = '  -2.73%  '
$ws.Range("D24").Value = '1.000'
$ws.Range("E24").Value = '  -0.23%  '
$ws.Range("D25").Value = '0.1574'
$ws.Range("E25").Value = '  -3.35%  '
$ws.Range("D26").Value = '9.474'
$ws.Range("E26").Value = '  +1.59%  '
$ws.Range("D27").Value = '162.48'
$ws.Range("E27").Value = '  -2.08%  '
$ws.Range("D28").Value = '18.77'
$ws.Range("E28").Value = '  +0.68%  '
$ws.Range("D29").Value = '2.032'
$ws.Range("E29").Value = '  -0.41%  '
$ws.Range("D30").Value = '1.425'
$ws.Range("E30").Value = '  +1.66%  '
$ws.Range("D31").Value = '1.544'
$ws.Range("E31").Value = '  +0.38%  '
$ws.Range("D32").Value = '4.514'
$ws.Range("E32").Value = '  +2.23%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").Value = '4.092'
$ws.Range("E33").Value = '  +1.61%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.05532'
$ws.Range("E34").Value = '  -1.46%  '
$ws.Range("D35").Value = '1.248'
$ws.Range("E35").Value = '  -0.80%  '
$ws.Range("D36").Value = '0.7503'
$ws.Range("E36").Value = '  +1.97%  '
$ws.Range("D37").Value = '0.9971'
$ws.Range("E37").Value = '  -0.20%  '
$ws.Range("D38").Value = '2.686'
$ws.Range("E38").Value = '  +1.85%  '
$ws.Range("D39").Value = '0.01929'
$ws.Range("E39").Value = '  +1.48%  '
$ws.Range("D40").Value = '2.790'
$ws.Range("E40").Value = '  +0.38%  '
$ws.Range("D41").Value = '0.4496'
$ws.Range("E41").Value = '  +2.60%  '
$ws.Range("D42").Value = '74.64'
$ws.Range("E42").Value = '  +3.73%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '6.032'
$ws.Range("E43").Value = '  +4.37%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '1.079.18'
$ws.Range("E44").Value = '  +3.93%  '
$ws.Range("D45").Value = '0.8521'
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("E46").Value = '  -0.14%  '
$ws.Range("D47").Value = '1.889'
$ws.Range("E47").Value = '  +1.61%  '
$ws.Range("D48").Value = '102.68'
$ws.Range("E48").Value = '  -0.35%  '
$ws.Range("D49").Value = '7.614'
$ws.Range("E49").Value = '  +3.12%  '
$ws.Range("D50").Value = '9.787'
$ws.Range("E50").Value = '  -1.68%  '
$ws.Range("D51").Value = '2.989'
$ws.Range("E51").Value = '  +4.07%  '
